$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First negative test with Data Driven Test approach:
# LastName becomes String.Empty, MartialStatus/Hobbies flags all become true.
$ws.Range("C2").Value = "String.Empty"
$ws.Range("D2").Value = "true true true"
$ws.Range("E2").Value = "true true true"

# Move the active selection to E2
$ws.Range("E2").Select()
